# Updates "想去人数" (F column) counts across the 展览/演出/全部类型 sheets
# to reflect the regenerated gh-pages data snapshot.
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 630
$wsExhibit.Range("F3").Value = 656
$wsExhibit.Range("F4").Value = 923
$wsExhibit.Range("F5").Value = 682
$wsExhibit.Range("F6").Value = 822
$wsExhibit.Range("F7").Value = 379
$wsExhibit.Range("F8").Value = 585
$wsExhibit.Range("F10").Value = 1180
$wsExhibit.Range("F11").Value = 611
$wsExhibit.Range("F12").Value = 362
$wsExhibit.Range("F13").Value = 488
$wsExhibit.Range("F14").Value = 161
$wsExhibit.Range("F15").Value = 245
$wsExhibit.Range("F16").Value = 325
$wsExhibit.Range("F18").Value = 79
$wsExhibit.Range("F19").Value = 540
$wsExhibit.Range("F20").Value = 54
$wsExhibit.Range("F21").Value = 551
$wsExhibit.Range("F22").Value = 22
$wsExhibit.Range("F23").Value = 594

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F9").Value = 212
$wsShow.Range("F13").Value = 53

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 630
$wsAll.Range("F7").Value = 656
$wsAll.Range("F8").Value = 923
$wsAll.Range("F9").Value = 682
$wsAll.Range("F10").Value = 822
$wsAll.Range("F11").Value = 379
$wsAll.Range("F12").Value = 585
$wsAll.Range("F14").Value = 1180
$wsAll.Range("F15").Value = 611
$wsAll.Range("F18").Value = 362
$wsAll.Range("F19").Value = 488
$wsAll.Range("F21").Value = 161
$wsAll.Range("F22").Value = 245
$wsAll.Range("F24").Value = 325
$wsAll.Range("F26").Value = 79
$wsAll.Range("F27").Value = 212
$wsAll.Range("F29").Value = 540
$wsAll.Range("F32").Value = 53
$wsAll.Range("F33").Value = 54
$wsAll.Range("F34").Value = 551
$wsAll.Range("F35").Value = 22
$wsAll.Range("F36").Value = 594
